# "Added button to unallocate a student from a project"
# Underlying data effect on the export sheet:
#  - Student names refreshed (sample/test data regenerated)
#  - The allocation that used to live on row 5 (Immanuel Wuckert) moved to
#    row 3 (now Ronny Dickinson): project columns C/D/E swapped between the
#    two rows, i.e. row 3 becomes allocated and row 5 becomes unallocated.
#  - Column A auto-fit width shrank slightly to match the new (shorter)
#    longest name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Aimee Collins"
$ws.Range("A3").Value = "Ronny Dickinson"
$ws.Range("A4").Value = "Al Herzog"
$ws.Range("A5").Value = "Loyal Steuber"

# Row 3 becomes allocated to the project
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# Row 5 becomes unallocated from the project
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Column A re-fits to the new longest name
$ws.Columns.Item(1).ColumnWidth = 18.709717
